# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E) on rows 16-24 lists the debt-period years
# for the worker. The previous account-statement periods are replaced with
# a new set of periods (the same nine years, now listed in the opposite/
# descending order), reflecting that the old EC (estado de cuenta) periods
# were removed and new ones added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2101"
$ws.Range("E17").Value = "2012"
$ws.Range("E18").Value = "2011"
$ws.Range("E19").Value = "2010"
$ws.Range("E20").Value = "2009"
$ws.Range("E21").Value = "2008"
$ws.Range("E22").Value = "2007"
$ws.Range("E23").Value = "2006"
$ws.Range("E24").Value = "2005"
